$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 19251
$ws.Range("B5").Value = "Pete Alonso"
$ws.Range("C5").Value = "NYM"
$ws.Range("D5").Value = 685
$ws.Range("E5").Value = 40
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 96
$ws.Range("H5").Value = 83
$ws.Range("I5").Value = 117
$ws.Range("J5").Value = 159
$ws.Range("K5").Value = 143
$ws.Range("L5").Value = 59
$ws.Range("M5").Value = 0.869
$ws.Range("N5").Value = 0.354
$ws.Range("O5").Value = 12.3
$ws.Range("P5").Value = 32
$ws.Range("Q5").Value = 116.5

$ws.Range("A6").Value = 17350
$ws.Range("B6").Value = "Rafael Devers"
$ws.Range("C6").Value = "BOS"
$ws.Range("D6").Value = 614
$ws.Range("E6").Value = 27
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 113
$ws.Range("H6").Value = 84
$ws.Range("I6").Value = 102
$ws.Range("J6").Value = 151
$ws.Range("K6").Value = 141
$ws.Range("L6").Value = 51
$ws.Range("M6").Value = 0.879
$ws.Range("N6").Value = 0.361
$ws.Range("O6").Value = 11.5
$ws.Range("P6").Value = 26.8
$ws.Range("Q6").Value = 113.7

$ws.Range("A7").Value = 18345
$ws.Range("B7").Value = "Kyle Tucker"
$ws.Range("C7").Value = "HOU"
$ws.Range("D7").Value = 609
$ws.Range("E7").Value = 30
$ws.Range("F7").Value = 25
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = 70
$ws.Range("I7").Value = 122
$ws.Range("J7").Value = 148
$ws.Range("K7").Value = 129
$ws.Range("L7").Value = 46
$ws.Range("M7").Value = 0.808
$ws.Range("N7").Value = 0.353
$ws.Range("O7").Value = 10.1
$ws.Range("P7").Value = 23.9
$ws.Range("Q7").Value = 111.3

$ws.Range("A8").Value = 14162
$ws.Range("B8").Value = "Carlos Correa"
$ws.Range("C8").Value = "MIN"
$ws.Range("D8").Value = 590
$ws.Range("E8").Value = 22
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 117
$ws.Range("H8").Value = 92
$ws.Range("I8").Value = 130
$ws.Range("J8").Value = 118
$ws.Range("K8").Value = 140
$ws.Range("L8").Value = 46
$ws.Range("M8").Value = 0.834
$ws.Range("N8").Value = 0.363
$ws.Range("O8").Value = 11.4
$ws.Range("P8").Value = 20.7
$ws.Range("Q8").Value = 114.6

$ws.Range("A9").Value = 19611
$ws.Range("B9").Value = "Vladimir Guerrero Jr."
$ws.Range("C9").Value = "TOR"
$ws.Range("D9").Value = 706
$ws.Range("E9").Value = 32
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = 74
$ws.Range("I9").Value = 103
$ws.Range("J9").Value = 138
$ws.Range("K9").Value = 132
$ws.Range("L9").Value = 59
$ws.Range("M9").Value = 0.818
$ws.Range("N9").Value = 0.348
$ws.Range("O9").Value = 11.2
$ws.Range("P9").Value = 19.3
$ws.Range("Q9").Value = 118.4

$ws.Range("A10").Value = 11579
$ws.Range("B10").Value = "Bryce Harper"
$ws.Range("C10").Value = "PHI"
$ws.Range("D10").Value = 426
$ws.Range("E10").Value = 18
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 111
$ws.Range("H10").Value = 90
$ws.Range("I10").Value = 129
$ws.Range("J10").Value = 146
$ws.Range("K10").Value = 138
$ws.Range("L10").Value = 37
$ws.Range("M10").Value = 0.877
$ws.Range("N10").Value = 0.38
$ws.Range("O10").Value = 12.8
$ws.Range("P10").Value = 18.1
$ws.Range("Q10").Value = 114.3

$ws.Range("A11").Value = 13419
$ws.Range("B11").Value = "Christian Walker"
$ws.Range("C11").Value = "ARI"
$ws.Range("D11").Value = 667
$ws.Range("E11").Value = 36
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 85
$ws.Range("H11").Value = 87
$ws.Range("I11").Value = 124
$ws.Range("J11").Value = 152
$ws.Range("K11").Value = 122
$ws.Range("L11").Value = 53
$ws.Range("M11").Value = 0.804
$ws.Range("N11").Value = 0.359
$ws.Range("O11").Value = 11.5
$ws.Range("P11").Value = 16.9
$ws.Range("Q11").Value = 112.6

$ws.Range("A12").Value = 3473
$ws.Range("B12").Value = "Anthony Rizzo"
$ws.Range("C12").Value = "NYY"
$ws.Range("D12").Value = 548
$ws.Range("E12").Value = 32
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 74
$ws.Range("H12").Value = 83
$ws.Range("I12").Value = 133
$ws.Range("J12").Value = 171
$ws.Range("K12").Value = 133
$ws.Range("L12").Value = 40
$ws.Range("M12").Value = 0.817
$ws.Range("N12").Value = 0.349
$ws.Range("O12").Value = 10.9
$ws.Range("P12").Value = 16.4
$ws.Range("Q12").Value = 113.3

$ws.Range("A13").Value = 11609
$ws.Range("B13").Value = "Willson Contreras"
$ws.Range("C13").Value = "CHC"
$ws.Range("D13").Value = 487
$ws.Range("E13").Value = 22
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 93
$ws.Range("H13").Value = 94
$ws.Range("I13").Value = 111
$ws.Range("J13").Value = 144
$ws.Range("K13").Value = 132
$ws.Range("L13").Value = 33
$ws.Range("M13").Value = 0.815
$ws.Range("N13").Value = 0.364
$ws.Range("O13").Value = 10.5
$ws.Range("P13").Value = 15.2
$ws.Range("Q13").Value = 116.2

$ws.Range("A14").Value = 19197
$ws.Range("B14").Value = "Will Smith"
$ws.Range("C14").Value = "LAD"
$ws.Range("D14").Value = 578
$ws.Range("E14").Value = 24
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 95
$ws.Range("H14").Value = 73
$ws.Range("I14").Value = 116
$ws.Range("J14").Value = 132
$ws.Range("K14").Value = 127
$ws.Range("L14").Value = 43
$ws.Range("M14").Value = 0.807
$ws.Range("N14").Value = 0.352
$ws.Range("O14").Value = 10.3
$ws.Range("P14").Value = 14.9
$ws.Range("Q14").Value = 109

$ws.Range("A15").Value = 16535
$ws.Range("B15").Value = "Danny Jansen"
$ws.Range("C15").Value = "TOR"
$ws.Range("D15").Value = 248
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 88
$ws.Range("H15").Value = 80
$ws.Range("I15").Value = 127
$ws.Range("J15").Value = 171
$ws.Range("K15").Value = 140
$ws.Range("L15").Value = 23
$ws.Range("M15").Value = 0.855
$ws.Range("N15").Value = 0.368
$ws.Range("O15").Value = 13.1
$ws.Range("P15").Value = 12
$ws.Range("Q15").Value = 109.3
